$wb = $excel.ActiveWorkbook

# Sheets that contain this dataset: "展览" (position 1) and "全部类型" (position 4)
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 8576
    $ws.Range("F4").Value = 1531
    $ws.Range("F7").Value = 272
    $ws.Range("F10").Value = 129
    $ws.Range("F12").Value = 463
    $ws.Range("F13").Value = 1275
    $ws.Range("F14").Value = 290
    $ws.Range("F15").Value = 85
    $ws.Range("F16").Value = 151
    $ws.Range("F20").Value = 126
}
